$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 32   Number  30"
$ws.Range("C9").Value = "Report Covering the Week  7/21/2025  Through  7/27/2025"

# --- Crime statistics data updates ---
# Simple numeric value changes (style unchanged)
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 0
$ws.Range("G16").Value = 4
$ws.Range("H16").Value = -50
$ws.Range("I16").Value = 17
$ws.Range("K16").Value = -48.484848484848
$ws.Range("L16").Value = -56.410256410256
$ws.Range("M16").Value = -65.306122448979
$ws.Range("N16").Value = -91.052631578947
$ws.Range("F17").Value = 9
$ws.Range("G17").Value = 6
$ws.Range("H17").Value = 50
$ws.Range("I17").Value = 71
$ws.Range("J17").Value = 57
$ws.Range("K17").Value = 24.561403508771
$ws.Range("L17").Value = 26.785714285714
$ws.Range("M17").Value = 136.666666666667
$ws.Range("N17").Value = 4.411764705882
$ws.Range("C18").Value = 7
$ws.Range("E18").Value = 40
$ws.Range("F18").Value = 18
$ws.Range("G18").Value = 18
$ws.Range("I18").Value = 165
$ws.Range("J18").Value = 146
$ws.Range("K18").Value = 13.013698630137
$ws.Range("L18").Value = -7.821229050279
$ws.Range("M18").Value = 25
$ws.Range("N18").Value = -71.938775510204
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 57.142857142857
$ws.Range("F19").Value = 33
$ws.Range("G19").Value = 34
$ws.Range("H19").Value = -2.941176470588
$ws.Range("I19").Value = 254
$ws.Range("J19").Value = 269
$ws.Range("K19").Value = -5.576208178438
$ws.Range("L19").Value = -31.903485254691
$ws.Range("M19").Value = 28.934010152284
$ws.Range("N19").Value = -15.050167224080
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 9
$ws.Range("E20").Value = -22.222222222222
$ws.Range("F20").Value = 22
$ws.Range("G20").Value = 25
$ws.Range("H20").Value = -12
$ws.Range("I20").Value = 153
$ws.Range("J20").Value = 149
$ws.Range("K20").Value = 2.684563758389
$ws.Range("L20").Value = 71.910112359550
$ws.Range("M20").Value = 91.25
$ws.Range("N20").Value = -91.662125340599
$ws.Range("C21").Value = 27
$ws.Range("E21").Value = 12.5
$ws.Range("F21").Value = 85
$ws.Range("G21").Value = 89
$ws.Range("H21").Value = -4.494382022471
$ws.Range("I21").Value = 666
$ws.Range("J21").Value = 660
$ws.Range("K21").Value = 0.909090909090
$ws.Range("L21").Value = -10.604026845637
$ws.Range("M21").Value = 35.365853658536
$ws.Range("N21").Value = -77.718300434928
$ws.Range("C24").Value = 14
$ws.Range("D24").Value = 9
$ws.Range("E24").Value = 55.555555555555
$ws.Range("F24").Value = 50
$ws.Range("G24").Value = 46
$ws.Range("H24").Value = 8.695652173913
$ws.Range("I24").Value = 321
$ws.Range("J24").Value = 323
$ws.Range("K24").Value = -0.619195046439
$ws.Range("L24").Value = -1.230769230769
$ws.Range("M24").Value = 21.132075471698
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 5
$ws.Range("G25").Value = 11
$ws.Range("H25").Value = -54.545454545454
$ws.Range("I25").Value = 55
$ws.Range("J25").Value = 67
$ws.Range("K25").Value = -17.910447761194
$ws.Range("L25").Value = -24.657534246575
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = -20
$ws.Range("F26").Value = 25
$ws.Range("G26").Value = 21
$ws.Range("H26").Value = 19.047619047619
$ws.Range("I26").Value = 118
$ws.Range("J26").Value = 122
$ws.Range("K26").Value = -3.278688524590
$ws.Range("L26").Value = -0.840336134453
$ws.Range("M26").Value = 21.649484536082
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 0

# Text -> Number conversions (set value + number format to match target style)
$ws.Range("D14").Value = 1
$ws.Range("D14").NumberFormat = "#,##0"
$ws.Range("E14").Value = -100
$ws.Range("E14").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G14").Value = 1
$ws.Range("G14").NumberFormat = "#,##0"
$ws.Range("H14").Value = -100
$ws.Range("H14").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("J14").Value = 1
$ws.Range("J14").NumberFormat = "#,##0"
$ws.Range("K14").Value = 100
$ws.Range("K14").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("C25").Value = 2
$ws.Range("C25").NumberFormat = "#,##0"
$ws.Range("D29").Value = 1
$ws.Range("D29").NumberFormat = "#,##0"
$ws.Range("E29").Value = -100
$ws.Range("E29").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G29").Value = 1
$ws.Range("G29").NumberFormat = "#,##0"
$ws.Range("H29").Value = -100
$ws.Range("H29").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("J29").Value = 1
$ws.Range("J29").NumberFormat = "#,##0"
$ws.Range("K29").Value = 100
$ws.Range("K29").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("D30").Value = 1
$ws.Range("D30").NumberFormat = "#,##0"
$ws.Range("E30").Value = -100
$ws.Range("E30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G30").Value = 1
$ws.Range("G30").NumberFormat = "#,##0"
$ws.Range("H30").Value = -100
$ws.Range("H30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("J30").Value = 1
$ws.Range("J30").NumberFormat = "#,##0"
$ws.Range("K30").Value = 100
$ws.Range("K30").NumberFormat = "#,##0.0;""-""#,##0.0"

# Number -> Text conversions (force text via "@" format, then restore style 13 via PasteSpecial formats)
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "0"
$ws.Range("C22").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0"
$ws.Range("C22").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "0"
$ws.Range("C22").Copy()
$ws.Range("G28").PasteSpecial(-4122)
$ws.Range("H28").NumberFormat = "@"
$ws.Range("H28").Value = "***.*"
$ws.Range("E22").Copy()
$ws.Range("H28").PasteSpecial(-4122)
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "0"
$ws.Range("C22").Copy()
$ws.Range("G33").PasteSpecial(-4122)
$ws.Range("H33").NumberFormat = "@"
$ws.Range("H33").Value = "***.*"
$ws.Range("E22").Copy()
$ws.Range("H33").PasteSpecial(-4122)
